# Weekly price update: insert two new rows (Primera/Segunda for "Crespo record")
# above row 117, pushing the existing rows 117:193 down to 119:195.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("117:118").Insert()

# New row 117 - Crespo record / Primera
$ws.Cells.Item(117, 1).Value2 = 11
$ws.Cells.Item(117, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(117, 3).Value2 = "Bíobío"
$ws.Cells.Item(117, 4).Value2 = 44455
$ws.Cells.Item(117, 5).Value2 = 8
$ws.Cells.Item(117, 6).Value2 = 100112006
$ws.Cells.Item(117, 7).Value2 = "Repollo"
$ws.Cells.Item(117, 8).Value2 = "Crespo record"
$ws.Cells.Item(117, 9).Value2 = "Primera"
$ws.Cells.Item(117, 10).Value2 = 1000
$ws.Cells.Item(117, 11).Value2 = 700
$ws.Cells.Item(117, 12).Value2 = 800
$ws.Cells.Item(117, 13).Value2 = 750
$ws.Cells.Item(117, 14).Value2 = "`$/unidad"
$ws.Cells.Item(117, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(117, 16).Value2 = 750
$ws.Cells.Item(117, 17).Value2 = 1
$ws.Cells.Item(117, 18).Value2 = "Hortaliza"

# New row 118 - Crespo record / Segunda
$ws.Cells.Item(118, 1).Value2 = 11
$ws.Cells.Item(118, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(118, 3).Value2 = "Bíobío"
$ws.Cells.Item(118, 4).Value2 = 44455
$ws.Cells.Item(118, 5).Value2 = 8
$ws.Cells.Item(118, 6).Value2 = 100112006
$ws.Cells.Item(118, 7).Value2 = "Repollo"
$ws.Cells.Item(118, 8).Value2 = "Crespo record"
$ws.Cells.Item(118, 9).Value2 = "Segunda"
$ws.Cells.Item(118, 10).Value2 = 500
$ws.Cells.Item(118, 11).Value2 = 600
$ws.Cells.Item(118, 12).Value2 = 600
$ws.Cells.Item(118, 13).Value2 = 600
$ws.Cells.Item(118, 14).Value2 = "`$/unidad"
$ws.Cells.Item(118, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(118, 16).Value2 = 600
$ws.Cells.Item(118, 17).Value2 = 1
$ws.Cells.Item(118, 18).Value2 = "Hortaliza"
